# Update the data table (A2:K17) with the new HOURS window (24..40) and
# append 5 new trailing rows (13..17) that didn't exist before.
#
# Columns: A, B, C, D, E, F, G, H, I, J, K
$data = @(
    @(24, 25, 165.1488800499209, 0, 0, 0, 247.7233200748814, 0.4, 0.48, 0.2, 0.4),
    @(25, 26, 130.2933375268219, 0, 0, 0, 195.4400062902328, 0.37, 0.46, 0.1, 0.33),
    @(26, 27, 103.6774555661557, 0, 0, 0, 155.5161833492336, 0.38, 0.48, 0.19, 0.32),
    @(27, 28, 97.78916096758275, 0, 0, 0, 146.6837414513741, 0.35, 0.43, 0.14, 0.32),
    @(28, 29, 94.78573553475324, 0, 0, 0, 142.1786033021299, 0.4, 0.49, 0.14, 0.39),
    @(29, 30, 98.71785172641819, 0, 0, 0, 148.0767775896273, 0.37, 0.46, 0.15, 0.33),
    @(30, 31, 130.7873219730109, 0, 0, 0, 196.1809829595164, 0.37, 0.49, 0.13, 0.4),
    @(31, 32, 201.8420447128433, 0, 0, 0, 302.7630670692649, 0.37, 0.46, 0.19, 0.31),
    @(32, 33, 294.8296768634683, 0.00277778, 0.00277778, 0.00555556, 442.2445152952025, 0.34, 0.46, 0.18, 0.38),
    @(33, 34, 349.0296502993302, 0.0138889, 0.0138889, 0.0277778, 523.5444754489954, 0.39, 0.5, 0.15, 0.4),
    @(34, 35, 364.1653337305623, 0.03055558, 0.03055558, 0.06111116, 546.2480005958435, 0.38, 0.43, 0.11, 0.37),
    @(35, 36, 381.7906987705871, 0.05833338, 0.05833338, 0.11666676, 572.6860481558806, 0.37, 0.49, 0.12, 0.3),
    @(36, 37, 410.3627591381613, 0.07777784, 0.07777784, 0.15555568, 615.5441387072419, 0.4, 0.46, 0.19, 0.33),
    @(37, 38, 408.4460994869478, 0.07222228, 0.07222228, 0.14444456, 612.6691492304217, 0.33, 0.45, 0.11, 0.32),
    @(38, 39, 377.9573794681603, 0.0555556, 0.0555556, 0.1111112, 566.9360692022404, 0.37, 0.48, 0.14, 0.33),
    @(39, 40, 354.1473291618483, 0.01111112, 0.01111112, 0.02222224, 531.2209937427724, 0.31, 0.44, 0.18, 0.32)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastOriginalRow = 12

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $firstRow + $i
    $rowVals = $data[$i]

    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $cell = $ws.Cells.Item($rowNum, $c + 1)
        $cell.Value = $rowVals[$c]
    }

    if ($rowNum -gt $lastOriginalRow) {
        # Brand-new row: column A needs the same style (bold, centered,
        # top-aligned, thin border) the existing A2:A12 cells already use.
        $aCell = $ws.Cells.Item($rowNum, 1)
        $aCell.Borders.LineStyle = 1
        $aCell.Font.Bold = $true
        $aCell.HorizontalAlignment = -4108
        $aCell.VerticalAlignment = -4160
    }
}
